# Insert a new weekly price record as row 77 ("Fruta / hortaliza, semanal"),
# pushing the existing rows 77-124 down to 78-125 (dimension grows to A1:R125).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 77..124 down by one, leaving a blank row 77 to populate.
$ws.Rows.Item(77).Insert()

# Populate the new row 77 with the new weekly record.
$ws.Range("A77").Value = 3
$ws.Range("B77").Value = "Femacal de La Calera"
$ws.Range("C77").Value = "Coquimbo"
$ws.Range("D77").Value = 44567
$ws.Range("E77").Value = 5
$ws.Range("F77").Value = 100112030
$ws.Range("G77").Value = "Poroto granado"
$ws.Range("H77").Value = "Sin especificar"
$ws.Range("I77").Value = "Primera"
$ws.Range("J77").Value = 80
$ws.Range("K77").Value = 22000
$ws.Range("L77").Value = 23000
$ws.Range("M77").Value = 22562
$ws.Range("N77").Value = "$/malla 25 kilos"
$ws.Range("O77").Value = "Provincia de Quillota"
$ws.Range("P77").Value = 902
$ws.Range("Q77").Value = 25
$ws.Range("R77").Value = "Hortaliza"
